$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new rows (bottom-to-top, using original pre-edit row positions) to make room for new weekly data
$ws.Range("A118:A119").EntireRow.Insert()
$ws.Range("A80:A83").EntireRow.Insert()
$ws.Range("A40:A41").EntireRow.Insert()

# Populate the newly inserted rows with their data
$ws.Range("A40").Value = 11
$ws.Range("B40").Value = "Vega Monumental Concepción"
$ws.Range("C40").Value = "Bíobío"
$ws.Range("D40").Value2 = 44434
$ws.Range("E40").Value = 8
$ws.Range("F40").Value = 100112009
$ws.Range("G40").Value = "Acelga"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 200
$ws.Range("K40").Value = 600
$ws.Range("L40").Value = 700
$ws.Range("M40").Value = 650
$ws.Range("N40").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O40").Value = "Región de Ñuble"
$ws.Range("P40").Value = 650
$ws.Range("Q40").Value = 1
$ws.Range("R40").Value = "Hortaliza"

$ws.Range("A41").Value = 11
$ws.Range("B41").Value = "Vega Monumental Concepción"
$ws.Range("C41").Value = "Bíobío"
$ws.Range("D41").Value2 = 44434
$ws.Range("E41").Value = 8
$ws.Range("F41").Value = 100112009
$ws.Range("G41").Value = "Acelga"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Segunda"
$ws.Range("J41").Value = 100
$ws.Range("K41").Value = 500
$ws.Range("L41").Value = 500
$ws.Range("M41").Value = 500
$ws.Range("N41").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O41").Value = "Región de Ñuble"
$ws.Range("P41").Value = 500
$ws.Range("Q41").Value = 1
$ws.Range("R41").Value = "Hortaliza"

$ws.Range("A82").Value = 11
$ws.Range("B82").Value = "Vega Monumental Concepción"
$ws.Range("C82").Value = "Bíobío"
$ws.Range("D82").Value2 = 44435
$ws.Range("E82").Value = 8
$ws.Range("F82").Value = 100112009
$ws.Range("G82").Value = "Acelga"
$ws.Range("H82").Value = "Sin especificar"
$ws.Range("I82").Value = "Primera"
$ws.Range("J82").Value = 200
$ws.Range("K82").Value = 600
$ws.Range("L82").Value = 700
$ws.Range("M82").Value = 650
$ws.Range("N82").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O82").Value = "Región Metropolitana"
$ws.Range("P82").Value = 650
$ws.Range("Q82").Value = 1
$ws.Range("R82").Value = "Hortaliza"

$ws.Range("A83").Value = 11
$ws.Range("B83").Value = "Vega Monumental Concepción"
$ws.Range("C83").Value = "Bíobío"
$ws.Range("D83").Value2 = 44435
$ws.Range("E83").Value = 8
$ws.Range("F83").Value = 100112009
$ws.Range("G83").Value = "Acelga"
$ws.Range("H83").Value = "Sin especificar"
$ws.Range("I83").Value = "Primera"
$ws.Range("J83").Value = 200
$ws.Range("K83").Value = 600
$ws.Range("L83").Value = 700
$ws.Range("M83").Value = 650
$ws.Range("N83").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O83").Value = "Región de Ñuble"
$ws.Range("P83").Value = 650
$ws.Range("Q83").Value = 1
$ws.Range("R83").Value = "Hortaliza"

$ws.Range("A84").Value = 11
$ws.Range("B84").Value = "Vega Monumental Concepción"
$ws.Range("C84").Value = "Bíobío"
$ws.Range("D84").Value2 = 44435
$ws.Range("E84").Value = 8
$ws.Range("F84").Value = 100112009
$ws.Range("G84").Value = "Acelga"
$ws.Range("H84").Value = "Sin especificar"
$ws.Range("I84").Value = "Segunda"
$ws.Range("J84").Value = 100
$ws.Range("K84").Value = 500
$ws.Range("L84").Value = 500
$ws.Range("M84").Value = 500
$ws.Range("N84").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O84").Value = "Región Metropolitana"
$ws.Range("P84").Value = 500
$ws.Range("Q84").Value = 1
$ws.Range("R84").Value = "Hortaliza"

$ws.Range("A85").Value = 11
$ws.Range("B85").Value = "Vega Monumental Concepción"
$ws.Range("C85").Value = "Bíobío"
$ws.Range("D85").Value2 = 44435
$ws.Range("E85").Value = 8
$ws.Range("F85").Value = 100112009
$ws.Range("G85").Value = "Acelga"
$ws.Range("H85").Value = "Sin especificar"
$ws.Range("I85").Value = "Segunda"
$ws.Range("J85").Value = 100
$ws.Range("K85").Value = 500
$ws.Range("L85").Value = 500
$ws.Range("M85").Value = 500
$ws.Range("N85").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O85").Value = "Región de Ñuble"
$ws.Range("P85").Value = 500
$ws.Range("Q85").Value = 1
$ws.Range("R85").Value = "Hortaliza"

$ws.Range("A124").Value = 11
$ws.Range("B124").Value = "Vega Monumental Concepción"
$ws.Range("C124").Value = "Bíobío"
$ws.Range("D124").Value2 = 44433
$ws.Range("E124").Value = 8
$ws.Range("F124").Value = 100112009
$ws.Range("G124").Value = "Acelga"
$ws.Range("H124").Value = "Sin especificar"
$ws.Range("I124").Value = "Primera"
$ws.Range("J124").Value = 200
$ws.Range("K124").Value = 600
$ws.Range("L124").Value = 700
$ws.Range("M124").Value = 650
$ws.Range("N124").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O124").Value = "Región Metropolitana"
$ws.Range("P124").Value = 650
$ws.Range("Q124").Value = 1
$ws.Range("R124").Value = "Hortaliza"

$ws.Range("A125").Value = 11
$ws.Range("B125").Value = "Vega Monumental Concepción"
$ws.Range("C125").Value = "Bíobío"
$ws.Range("D125").Value2 = 44433
$ws.Range("E125").Value = 8
$ws.Range("F125").Value = 100112009
$ws.Range("G125").Value = "Acelga"
$ws.Range("H125").Value = "Sin especificar"
$ws.Range("I125").Value = "Segunda"
$ws.Range("J125").Value = 100
$ws.Range("K125").Value = 500
$ws.Range("L125").Value = 500
$ws.Range("M125").Value = 500
$ws.Range("N125").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O125").Value = "Región Metropolitana"
$ws.Range("P125").Value = 500
$ws.Range("Q125").Value = 1
$ws.Range("R125").Value = "Hortaliza"
